# Update list of wireless networks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "WLAN_XXXX & JAZZTEL_XXXX / Comtrend/Zyxel?" row — the
# rows below it shuffle up to fill the gap.
$ws.Rows(38).Delete()

# Row 37 (WLAN_XXXX & JAZZTEL_XXXX / Comtrend) now supports the password
# algorithm too; drop the stale Encryption note and record it.
$ws.Range("E37").Value = ""
$ws.Range("F37").Value = "Yes"
$ws.Range("H37").Value = "Zyxel P660HW-D1"

# Two more Comtrend routers (moved up into rows 41/42) are now confirmed
# supported.
$ws.Range("F41").Value = "Yes"
$ws.Range("F42").Value = "Yes"

# The last row (formerly row 49) becomes a brand-new network entry.
$ws.Range("A48").Value = "WLAN_XX"
$ws.Range("B48").Value = "00:23:F8:XX:XX:XX"
$ws.Range("C48").Value = "?"
$ws.Range("H48").Value = "More than 1 password…"

# B48 is highlighted with a black-on-white fill to call it out.
$ws.Range("B48").Font.Color = 0
$ws.Range("B48").Interior.Color = 16777215
$ws.Range("B48").Interior.PatternColor = 0

# Restore the workbook / sheet view state captured in the saved file.
$excel.Width = 24920
$excel.Height = 12980
$excel.Left = 400
$excel.Top = 1380

$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
$ws.Range("B5").Select()
